$d = $word.ActiveDocument

$d.Content.Find.Execute("Kevin Chen (s3895923)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Kevin Chen (s3780646)", 2)
